# Add a new "canonical SMILES" column (D) to the microstate list sheet.
# Column D mirrors column C ("canonical isomeric SMILES") with any
# stereochemistry slash markers ("/" and "\") stripped out.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header
$ws.Range("D2").Value = "canonical SMILES"

# Data rows: 3..13 hold the per-microstate SMILES in column C.
for ($r = 3; $r -le 13; $r++) {
    $iso = $ws.Cells.Item($r, 3).Value2
    if ($iso -ne $null -and $iso -ne "") {
        $canonical = $iso -replace "/", "" -replace "\\", ""
        $ws.Cells.Item($r, 4).Value = $canonical
    }
}

# Match the new column's width to the authored workbook (closest value this
# host's column-width quantization can reach to the authored 36.85546875).
$ws.Columns.Item(4).ColumnWidth = 36
